$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 changes: Rent/300/45901.16680555556 -> Video Games/200/45924.16680555556
$ws.Range("A6").Value = "Video Games"
$ws.Range("B6").Value = 200
$ws.Range("C6").Value = 45924.16680555556

# Row 7 keeps Rent/300 but date now matches the old row6 date
$ws.Range("A7").Value = "Rent"
$ws.Range("B7").Value = 300
$ws.Range("C7").Value = 45901.16680555556

# New row 8: the old row7 data (Rent/300/45689.16680555556)
$ws.Range("A8").Value = "Rent"
$ws.Range("B8").Value = 300
$ws.Range("C8").Value = 45689.16680555556

# Match date number format (style) used by the other date cells in column C
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)
